$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 253, shifting existing rows 253-285 down to 254-286.
$ws.Rows("253:253").Insert()

# Populate the newly inserted row 253 with its values.
$ws.Range("A253").Value = 8
$ws.Range("B253").Value = "Terminal La Palmera de La Serena"
$ws.Range("C253").Value = "Coquimbo"
$ws.Range("D253").Value = 44748
$ws.Range("D253").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E253").Value = 4
$ws.Range("F253").Value = 100112012
$ws.Range("G253").Value = "Espinaca"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 3000
$ws.Range("K253").Value = 500
$ws.Range("L253").Value = 600
$ws.Range("M253").Value = 550
$ws.Range("N253").Value = "$/atado 300 a 500 gramos"
$ws.Range("O253").Value = "Provincia del Elquí"
$ws.Range("P253").Value = 1100
$ws.Range("Q253").Value = 0.5
$ws.Range("R253").Value = "Hortaliza"
